$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 98.912777
$ws.Range("H2").Value = 296.738331
$ws.Range("I2").Value = 0.8120825131376513
$ws.Range("J2").Value = 0.8120825131376513
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.922188
$ws.Range("N2").Value = 8.766564
$ws.Range("O2").Value = 0.04895158192025818
$ws.Range("P2").Value = 0.04895158192025818
$ws.Range("Q2").Value = 289.041729996076
$ws.Range("R2").Value = 2601.375569964684
$ws.Range("S2").Value = 0.03975272366786688
$ws.Range("T2").Value = 0.03975272366786688

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 98.912777
$ws.Range("H3").Value = 296.738331
$ws.Range("I3").Value = 0.8120825131376513
$ws.Range("J3").Value = 0.8120825131376513
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 32.47988333333333
$ws.Range("N3").Value = 97.43965
$ws.Range("O3").Value = 0.5440928748431295
$ws.Range("P3").Value = 0.5440928748431295
$ws.Range("Q3").Value = 3212.675457136017
$ws.Range("R3").Value = 28914.07911422415
$ws.Range("S3").Value = 0.4418483091828982
$ws.Range("T3").Value = 0.4418483091828982

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 98.912777
$ws.Range("H4").Value = 296.738331
$ws.Range("I4").Value = 0.8120825131376513
$ws.Range("J4").Value = 0.8120825131376513
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 24.293405
$ws.Range("N4").Value = 72.88021499999999
$ws.Range("O4").Value = 0.4069555432366123
$ws.Range("P4").Value = 0.4069555432366123
$ws.Range("Q4").Value = 2402.928151335685
$ws.Range("R4").Value = 21626.35336202116
$ws.Range("S4").Value = 0.3304814802868862
$ws.Range("T4").Value = 0.3304814802868862

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 17.04862266666667
$ws.Range("H5").Value = 51.14586800000001
$ws.Range("I5").Value = 0.1399706767982279
$ws.Range("J5").Value = 0.1399706767982279
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.922188
$ws.Range("N5").Value = 8.766564
$ws.Range("O5").Value = 0.04895158192025818
$ws.Range("P5").Value = 0.04895158192025818
$ws.Range("Q5").Value = 49.81928057306135
$ws.Range("R5").Value = 448.3735251575521
$ws.Range("S5").Value = 0.006851786051722437
$ws.Range("T5").Value = 0.006851786051722437

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 17.04862266666667
$ws.Range("H6").Value = 51.14586800000001
$ws.Range("I6").Value = 0.1399706767982279
$ws.Range("J6").Value = 0.1399706767982279
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 32.47988333333333
$ws.Range("N6").Value = 97.43965
$ws.Range("O6").Value = 0.5440928748431295
$ws.Range("P6").Value = 0.5440928748431295
$ws.Range("Q6").Value = 553.7372752073557
$ws.Range("R6").Value = 4983.6354768662
$ws.Range("S6").Value = 0.07615704793288637
$ws.Range("T6").Value = 0.07615704793288637

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 17.04862266666667
$ws.Range("H7").Value = 51.14586800000001
$ws.Range("I7").Value = 0.1399706767982279
$ws.Range("J7").Value = 0.1399706767982279
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 24.293405
$ws.Range("N7").Value = 72.88021499999999
$ws.Range("O7").Value = 0.4069555432366123
$ws.Range("P7").Value = 0.4069555432366123
$ws.Range("Q7").Value = 414.1690951335133
$ws.Range("R7").Value = 3727.52185620162
$ws.Range("S7").Value = 0.05696184281361914
$ws.Range("T7").Value = 0.05696184281361914

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.839988000000001
$ws.Range("H8").Value = 17.519964
$ws.Range("I8").Value = 0.0479468100641207
$ws.Range("J8").Value = 0.04794681006412069
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.922188
$ws.Range("N8").Value = 8.766564
$ws.Range("O8").Value = 0.04895158192025818
$ws.Range("P8").Value = 0.04895158192025818
$ws.Range("Q8").Value = 17.065542853744
$ws.Range("R8").Value = 153.589885683696
$ws.Range("S8").Value = 0.002347072200668864
$ws.Range("T8").Value = 0.002347072200668864

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.839988000000001
$ws.Range("H9").Value = 17.519964
$ws.Range("I9").Value = 0.0479468100641207
$ws.Range("J9").Value = 0.04794681006412069
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 32.47988333333333
$ws.Range("N9").Value = 97.43965
$ws.Range("O9").Value = 0.5440928748431295
$ws.Range("P9").Value = 0.5440928748431295
$ws.Range("Q9").Value = 189.6821289080667
$ws.Range("R9").Value = 1707.1391601726
$ws.Range("S9").Value = 0.02608751772734493
$ws.Range("T9").Value = 0.02608751772734492

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 5.839988000000001
$ws.Range("H10").Value = 17.519964
$ws.Range("I10").Value = 0.0479468100641207
$ws.Range("J10").Value = 0.04794681006412069
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 24.293405
$ws.Range("N10").Value = 72.88021499999999
$ws.Range("O10").Value = 0.4069555432366123
$ws.Range("P10").Value = 0.4069555432366123
$ws.Range("Q10").Value = 141.87319367914
$ws.Range("R10").Value = 1276.85874311226
$ws.Range("S10").Value = 0.01951222013610691
$ws.Range("T10").Value = 0.0195122201361069
